# Applies the diff: expands the "Enter your Raspberry Pi's private IP..."
# section with the new npm / node-red install walkthrough (headings,
# paragraphs, and sample-script text boxes), and relocates the _GoBack
# bookmark to wrap the final "node-red-admin hash-pw" text box instead of
# the old ", raspberry." run. The two trailing empty paragraphs at the end
# of the document are also dropped.

$d = $word.ActiveDocument

# Find the paragraph that starts the block we are rewriting.
$found = $d.Content
$found.Find.Execute("Enter your Raspberry Pi's private IP address into VNC Viewer")
$targetPara = $found.Paragraphs.Item(1)

# Replace everything from the start of that paragraph through the end of the
# document (covers the password paragraph and the two trailing blank ones)
# with the new block of paragraphs in one shot.
$rangeStart = $targetPara.Range.Start
$rangeEnd = $d.Content.End
$replaceRange = $d.Range($rangeStart, $rangeEnd)

$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><w:body><w:p><w:r><w:t>Enter your Raspberry Pi's private IP address into VNC Viewer.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">The default username and password are </w:t></w:r><w:r><w:t>pi</w:t></w:r><w:r><w:t>, raspberry.</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Install </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>npm</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">All new Raspbian OS install node-red, node.js by default. You just install NPM by yourself, for install additional node to node-red </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>api</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t>Run following command:</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:noProof/></w:rPr><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="06B88406" wp14:editId="6FAD8AB1"><wp:extent cx="5924550" cy="1404620"/><wp:effectExtent l="0" t="0" r="19050" b="20320"/><wp:docPr id="4" name="Text Box 2"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr txBox="1"><a:spLocks noChangeArrowheads="1"/></wps:cNvSpPr><wps:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="5924550" cy="1404620"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:solidFill><a:srgbClr val="FFFFFF"/></a:solidFill><a:ln w="9525"><a:solidFill><a:srgbClr val="000000"/></a:solidFill><a:miter lim="800000"/><a:headEnd/><a:tailEnd/></a:ln></wps:spPr><wps:txbx><w:txbxContent><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>sudo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> apt-get install </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>nodejs</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>npm</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:txbxContent></wps:txbx><wps:bodyPr rot="0" vert="horz" wrap="square" lIns="91440" tIns="45720" rIns="91440" bIns="45720" anchor="t" anchorCtr="0"><a:spAutoFit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic></wp:inline></w:drawing></mc:Choice><mc:Fallback><w:pict><v:shape w14:anchorId="06B88406" id="_x0000_s1028" type="#_x0000_t202" style="width:466.5pt;height:110.6pt;visibility:visible;mso-wrap-style:square;mso-left-percent:-10001;mso-top-percent:-10001;mso-position-horizontal:absolute;mso-position-horizontal-relative:char;mso-position-vertical:absolute;mso-position-vertical-relative:line;mso-left-percent:-10001;mso-top-percent:-10001;v-text-anchor:top"><v:textbox style="mso-fit-shape-to-text:t"><w:txbxContent><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>sudo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> apt-get install </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>nodejs</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>npm</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:txbxContent></v:textbox><w10:anchorlock/></v:shape></w:pict></mc:Fallback></mc:AlternateContent></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Running on Node-Red.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:noProof/></w:rPr><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="5B2CC3AC" wp14:editId="5F5B863C"><wp:extent cx="5924550" cy="1404620"/><wp:effectExtent l="0" t="0" r="19050" b="20320"/><wp:docPr id="5" name="Text Box 2"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr txBox="1"><a:spLocks noChangeArrowheads="1"/></wps:cNvSpPr><wps:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="5924550" cy="1404620"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:solidFill><a:srgbClr val="FFFFFF"/></a:solidFill><a:ln w="9525"><a:solidFill><a:srgbClr val="000000"/></a:solidFill><a:miter lim="800000"/><a:headEnd/><a:tailEnd/></a:ln></wps:spPr><wps:txbx><w:txbxContent><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>sudo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> node-red-start</w:t></w:r><w:r><w:tab/><w:t>//Start node-red</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>sudo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> node-red-stop</w:t></w:r><w:r><w:tab/><w:t>//stop node-red</w:t></w:r></w:p></w:txbxContent></wps:txbx><wps:bodyPr rot="0" vert="horz" wrap="square" lIns="91440" tIns="45720" rIns="91440" bIns="45720" anchor="t" anchorCtr="0"><a:spAutoFit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic></wp:inline></w:drawing></mc:Choice><mc:Fallback><w:pict><v:shape w14:anchorId="5B2CC3AC" id="_x0000_s1029" type="#_x0000_t202" style="width:466.5pt;height:110.6pt;visibility:visible;mso-wrap-style:square;mso-left-percent:-10001;mso-top-percent:-10001;mso-position-horizontal:absolute;mso-position-horizontal-relative:char;mso-position-vertical:absolute;mso-position-vertical-relative:line;mso-left-percent:-10001;mso-top-percent:-10001;v-text-anchor:top"><v:textbox style="mso-fit-shape-to-text:t"><w:txbxContent><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>sudo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> node-red-start</w:t></w:r><w:r><w:tab/><w:t>//Start node-red</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>sudo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> node-red-stop</w:t></w:r><w:r><w:tab/><w:t>//stop node-red</w:t></w:r></w:p></w:txbxContent></v:textbox><w10:anchorlock/></v:shape></w:pict></mc:Fallback></mc:AlternateContent></w:r></w:p><w:p><w:r><w:t>Access node-red API</w:t></w:r></w:p><w:p><w:r><w:rPr><w:noProof/></w:rPr><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="1BFB3DD5" wp14:editId="107FD3B2"><wp:extent cx="5924550" cy="1404620"/><wp:effectExtent l="0" t="0" r="19050" b="20320"/><wp:docPr id="6" name="Text Box 2"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr txBox="1"><a:spLocks noChangeArrowheads="1"/></wps:cNvSpPr><wps:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="5924550" cy="1404620"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:solidFill><a:srgbClr val="FFFFFF"/></a:solidFill><a:ln w="9525"><a:solidFill><a:srgbClr val="000000"/></a:solidFill><a:miter lim="800000"/><a:headEnd/><a:tailEnd/></a:ln></wps:spPr><wps:txbx><w:txbxContent><w:p><w:r><w:t xml:space="preserve">&lt;Raspberry pi </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ip</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> address&gt;:1880</w:t></w:r></w:p></w:txbxContent></wps:txbx><wps:bodyPr rot="0" vert="horz" wrap="square" lIns="91440" tIns="45720" rIns="91440" bIns="45720" anchor="t" anchorCtr="0"><a:spAutoFit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic></wp:inline></w:drawing></mc:Choice><mc:Fallback><w:pict><v:shape w14:anchorId="1BFB3DD5" id="_x0000_s1030" type="#_x0000_t202" style="width:466.5pt;height:110.6pt;visibility:visible;mso-wrap-style:square;mso-left-percent:-10001;mso-top-percent:-10001;mso-position-horizontal:absolute;mso-position-horizontal-relative:char;mso-position-vertical:absolute;mso-position-vertical-relative:line;mso-left-percent:-10001;mso-top-percent:-10001;v-text-anchor:top"><v:textbox style="mso-fit-shape-to-text:t"><w:txbxContent><w:p><w:r><w:t xml:space="preserve">&lt;Raspberry pi </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ip</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> address&gt;:1880</w:t></w:r></w:p></w:txbxContent></v:textbox><w10:anchorlock/></v:shape></w:pict></mc:Fallback></mc:AlternateContent></w:r></w:p><w:p><w:r><w:t>Editor &amp; Admin API security</w:t></w:r></w:p><w:p><w:r><w:t>The Editor and Admin API supports two types of authentication:</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">username/password </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>credential based</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> authentication</w:t></w:r></w:p><w:p><w:r><w:t>since Node-RED 0.17: authentication against any OAuth/OpenID provider such as Twitter or GitHub</w:t></w:r></w:p><w:p><w:r><w:t>Username/</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>password based</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> authentication</w:t></w:r></w:p><w:p><w:r><w:t>To enable user authentication on the Editor and Admin API, add the following to your settings.js file:</w:t></w:r><w:r><w:t xml:space="preserve"> (cd $HOME/pi</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>/.node</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>-red)</w:t></w:r></w:p><w:p><w:r><w:rPr><w:noProof/></w:rPr><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="7AC34878" wp14:editId="30C02825"><wp:extent cx="5924550" cy="1404620"/><wp:effectExtent l="0" t="0" r="19050" b="20320"/><wp:docPr id="7" name="Text Box 2"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr txBox="1"><a:spLocks noChangeArrowheads="1"/></wps:cNvSpPr><wps:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="5924550" cy="1404620"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:solidFill><a:srgbClr val="FFFFFF"/></a:solidFill><a:ln w="9525"><a:solidFill><a:srgbClr val="000000"/></a:solidFill><a:miter lim="800000"/><a:headEnd/><a:tailEnd/></a:ln></wps:spPr><wps:txbx><w:txbxContent><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>adminAuth</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: {</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    type: "credentials",</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    users: [{</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">        username: "admin",</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">        password: "$2a$08$zZWtXTja0fB1pzD4sHCMyOCMYz2Z6dNbM6tl8sJogENOMcxWV9DN.",</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">        permissions: "*"</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    }]</w:t></w:r></w:p><w:p><w:r><w:t>}</w:t></w:r></w:p></w:txbxContent></wps:txbx><wps:bodyPr rot="0" vert="horz" wrap="square" lIns="91440" tIns="45720" rIns="91440" bIns="45720" anchor="t" anchorCtr="0"><a:spAutoFit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic></wp:inline></w:drawing></mc:Choice><mc:Fallback><w:pict><v:shape w14:anchorId="7AC34878" id="_x0000_s1031" type="#_x0000_t202" style="width:466.5pt;height:110.6pt;visibility:visible;mso-wrap-style:square;mso-left-percent:-10001;mso-top-percent:-10001;mso-position-horizontal:absolute;mso-position-horizontal-relative:char;mso-position-vertical:absolute;mso-position-vertical-relative:line;mso-left-percent:-10001;mso-top-percent:-10001;v-text-anchor:top"><v:textbox style="mso-fit-shape-to-text:t"><w:txbxContent><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>adminAuth</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: {</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    type: "credentials",</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    users: [{</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">        username: "admin",</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">        password: "$2a$08$zZWtXTja0fB1pzD4sHCMyOCMYz2Z6dNbM6tl8sJogENOMcxWV9DN.",</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">        permissions: "*"</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    }]</w:t></w:r></w:p><w:p><w:r><w:t>}</w:t></w:r></w:p></w:txbxContent></v:textbox><w10:anchorlock/></v:shape></w:pict></mc:Fallback></mc:AlternateContent></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t>Username: admin, Password: password</w:t></w:r></w:p><w:p><w:r><w:t>Generating the password hash</w:t></w:r></w:p><w:p><w:r><w:t>To generate a suitable password hash, you can use the node-red-admin command-line tool:</w:t></w:r></w:p><w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:r><w:rPr><w:noProof/></w:rPr><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="47BA5D8A" wp14:editId="53A36414"><wp:extent cx="5924550" cy="1404620"/><wp:effectExtent l="0" t="0" r="19050" b="20320"/><wp:docPr id="8" name="Text Box 2"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr txBox="1"><a:spLocks noChangeArrowheads="1"/></wps:cNvSpPr><wps:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="5924550" cy="1404620"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:solidFill><a:srgbClr val="FFFFFF"/></a:solidFill><a:ln w="9525"><a:solidFill><a:srgbClr val="000000"/></a:solidFill><a:miter lim="800000"/><a:headEnd/><a:tailEnd/></a:ln></wps:spPr><wps:txbx><w:txbxContent><w:p><w:r><w:t>node-red-admin hash-pw</w:t></w:r></w:p></w:txbxContent></wps:txbx><wps:bodyPr rot="0" vert="horz" wrap="square" lIns="91440" tIns="45720" rIns="91440" bIns="45720" anchor="t" anchorCtr="0"><a:spAutoFit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic></wp:inline></w:drawing></mc:Choice><mc:Fallback><w:pict><v:shape w14:anchorId="47BA5D8A" id="_x0000_s1032" type="#_x0000_t202" style="width:466.5pt;height:110.6pt;visibility:visible;mso-wrap-style:square;mso-left-percent:-10001;mso-top-percent:-10001;mso-position-horizontal:absolute;mso-position-horizontal-relative:char;mso-position-vertical:absolute;mso-position-vertical-relative:line;mso-left-percent:-10001;mso-top-percent:-10001;v-text-anchor:top"><v:textbox style="mso-fit-shape-to-text:t"><w:txbxContent><w:p><w:r><w:t>node-red-admin hash-pw</w:t></w:r></w:p></w:txbxContent></v:textbox><w10:anchorlock/></v:shape></w:pict></mc:Fallback></mc:AlternateContent></w:r><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$replaceRange.InsertXML($xml)
